$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping Tag Glossary")

# Insert a new blank row at row 48 (pushes existing rows 48-54 down to 49-55)
$ws.Rows.Item(48).Insert()

# Copy the formatting (borders/fonts/alignment) from the row above (row 47,
# a normal interior data row) into the new row 48 so it matches the other
# data rows instead of the blank/default formatting left by Insert().
$ws.Range("A47:C47").Copy()
$ws.Range("A48:C48").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row with the new glossary entry.
$ws.Range("A48").Value = 46
$ws.Range("B48").Value = "Map:SDTM IGv3.2"
$ws.Range("C48").Value = "CDISC Study Data Tabulation Model Implementation Guide version 3.2"

# Renumber the "#" column for every row pushed down by the insert so the
# sequence stays contiguous (47, 48, 49, ... 53).
for ($r = 49; $r -le 55; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
